# Update gh-pages to output generated at 456a3b4
# This script updates the "想去人数" (F column) figures and one cover-image
# URL across the "展览" (sheet 1), "演出" (sheet 2) and "全部类型" (sheet 4)
# worksheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 275
$ws.Range("F5").Value = 272
$ws.Range("F6").Value = 1077
$ws.Range("F7").Value = 1410
$ws.Range("F9").Value = 101
$ws.Range("F12").Value = 134
$ws.Range("F13").Value = 123
$ws.Range("F15").Value = 1316
$ws.Range("F16").Value = 100
$ws.Range("F17").Value = 87
$ws.Range("F18").Value = 268
$ws.Range("F20").Value = 640
$ws.Range("F22").Value = 199
$ws.Range("F23").Value = 12
$ws.Range("F24").Value = 5648
$ws.Range("F25").Value = 54
$ws.Range("F29").Value = 14208
$ws.Range("F30").Value = 1417
$ws.Range("F31").Value = 194
$ws.Range("F32").Value = 92
$ws.Range("F33").Value = 85
$ws.Range("F34").Value = 425
$ws.Range("I34").Value = "//i0.hdslb.com/bfs/openplatform/202407/yw21E7Vn1721701909995.jpeg"
$ws.Range("F35").Value = 583
$ws.Range("F36").Value = 4174
$ws.Range("F37").Value = 119
$ws.Range("F38").Value = 353

# ---- Sheet 2: 演出 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 341

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 275
$ws.Range("F5").Value = 272
$ws.Range("F6").Value = 1077
$ws.Range("F7").Value = 1410
$ws.Range("F8").Value = 579
$ws.Range("F9").Value = 101
$ws.Range("F12").Value = 134
$ws.Range("F13").Value = 123
$ws.Range("F15").Value = 1316
$ws.Range("F16").Value = 100
$ws.Range("F17").Value = 87
$ws.Range("F18").Value = 268
$ws.Range("F19").Value = 341
$ws.Range("F21").Value = 640
$ws.Range("F24").Value = 199
$ws.Range("F25").Value = 12
$ws.Range("F27").Value = 5648
$ws.Range("F28").Value = 54
$ws.Range("F32").Value = 14209
$ws.Range("F33").Value = 1417
$ws.Range("F34").Value = 194
$ws.Range("F35").Value = 92
$ws.Range("F36").Value = 85
$ws.Range("F37").Value = 425
$ws.Range("I37").Value = "//i0.hdslb.com/bfs/openplatform/202407/yw21E7Vn1721701909995.jpeg"
$ws.Range("F38").Value = 583
$ws.Range("F39").Value = 4174
$ws.Range("F40").Value = 119
$ws.Range("F41").Value = 353
